# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections described by the diff to sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 21811.26
$ws.Range("I64").Value = 3247.2942
$ws.Range("J64").Value = 53370
$ws.Range("K64").Value = 3247.2942
$ws.Range("L64").Value = 53370
$ws.Range("M64").Value = -2999.2942
$ws.Range("N64").Value = -53866
$ws.Range("H67").Value = 21811.26
$ws.Range("I67").Value = 3247.2942
$ws.Range("J67").Value = 53370
$ws.Range("K67").Value = 3247.2942
$ws.Range("L67").Value = 53370
$ws.Range("M67").Value = -2389.2942
$ws.Range("N67").Value = -55086
$ws.Range("H74").Value = 6670034
$ws.Range("I74").Value = 7146222
$ws.Range("J74").Value = 3400
$ws.Range("K74").Value = 7146222
$ws.Range("L74").Value = 3400
$ws.Range("M74").Value = -7145286
$ws.Range("N74").Value = -5272
$ws.Range("H76").Value = 166671600
$ws.Range("I76").Value = 250002380
$ws.Range("J76").Value = 10004
$ws.Range("K76").Value = 250002380
$ws.Range("L76").Value = 10004
$ws.Range("M76").Value = -250002065
$ws.Range("N76").Value = -10634
$ws.Range("H77").Value = 6670034
$ws.Range("I77").Value = 7146222
$ws.Range("J77").Value = 3400
$ws.Range("K77").Value = 35731110
$ws.Range("L77").Value = 17000
$ws.Range("M77").Value = -35726430
$ws.Range("N77").Value = -26360
$ws.Range("H79").Value = 166671600
$ws.Range("I79").Value = 250002380
$ws.Range("J79").Value = 10004
$ws.Range("K79").Value = 250002380
$ws.Range("L79").Value = 10004
$ws.Range("M79").Value = -250001288
$ws.Range("N79").Value = -12188
$ws.Range("H100").Value = 11333
$ws.Range("I100").Value = 3999
$ws.Range("J100").Value = 15000
$ws.Range("K100").Value = 3999
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = -3458
$ws.Range("N100").Value = -16082
$ws.Range("H132").Value = 177796.16
$ws.Range("I132").Value = 2230.9773
$ws.Range("J132").Value = 772016.75
$ws.Range("K132").Value = 6692.9319
$ws.Range("L132").Value = 2316050.25
$ws.Range("M132").Value = -4162.9319
$ws.Range("N132").Value = -2321110.25
$ws.Range("H137").Value = 27342.578
$ws.Range("I137").Value = 36645.645
$ws.Range("K137").Value = 109936.935
$ws.Range("M137").Value = -107386.935
$ws.Range("H138").Value = 11889.64
$ws.Range("I138").Value = 714.3939
$ws.Range("J138").Value = 17651.875
$ws.Range("K138").Value = 2143.1817
$ws.Range("L138").Value = 52955.625
$ws.Range("M138").Value = 2996.8183
$ws.Range("N138").Value = -63235.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19341.031
$ws.Range("I32").Value = 3805.5686
$ws.Range("J32").Value = 80287.84
$ws.Range("K32").Value = 3805.5686
$ws.Range("L32").Value = 80287.84
$ws.Range("M32").Value = -3518.5686
$ws.Range("N32").Value = -80861.84
$ws.Range("H109").Value = 34987
$ws.Range("J109").Value = 34987
$ws.Range("L109").Value = 34987
$ws.Range("N109").Value = -37761
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 25500
$ws.Range("J108").Value = 25500
$ws.Range("L108").Value = 25500
$ws.Range("N108").Value = -33180
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H31").Value = 42174.19
$ws.Range("I31").Value = 72842.92999999999
$ws.Range("J31").Value = 6394
$ws.Range("K31").Value = 72842.92999999999
$ws.Range("L31").Value = 6394
$ws.Range("M31").Value = -72547.92999999999
$ws.Range("N31").Value = -6984
$ws.Range("H34").Value = 42174.19
$ws.Range("I34").Value = 72842.92999999999
$ws.Range("J34").Value = 6394
$ws.Range("K34").Value = 72842.92999999999
$ws.Range("L34").Value = 6394
$ws.Range("M34").Value = -72640.92999999999
$ws.Range("N34").Value = -6798
$ws.Range("H41").Value = 10152.167
$ws.Range("I41").Value = 5400
$ws.Range("J41").Value = 12528.25
$ws.Range("K41").Value = 5400
$ws.Range("L41").Value = 12528.25
$ws.Range("M41").Value = -4972
$ws.Range("N41").Value = -13384.25
$ws.Range("H50").Value = 10358.25
$ws.Range("J50").Value = 10358.25
$ws.Range("L50").Value = 10358.25
$ws.Range("N50").Value = -11608.25
$ws.Range("H51").Value = 10402.728
$ws.Range("J51").Value = 10402.728
$ws.Range("L51").Value = 10402.728
$ws.Range("N51").Value = -11874.728
$ws.Range("H52").Value = 77040
$ws.Range("J52").Value = 77040
$ws.Range("L52").Value = 77040
$ws.Range("N52").Value = -77628
$ws.Range("H58").Value = 3300.1904
$ws.Range("I58").Value = 771.7143
$ws.Range("J58").Value = 8357.143
$ws.Range("K58").Value = 771.7143
$ws.Range("L58").Value = 8357.143
$ws.Range("M58").Value = -568.7143
$ws.Range("N58").Value = -8763.143
$ws.Range("H60").Value = 9622.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 9622.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 9622.5
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -10644.5
$ws.Range("H61").Value = 10402.728
$ws.Range("J61").Value = 10402.728
$ws.Range("L61").Value = 10402.728
$ws.Range("N61").Value = -11098.728
$ws.Range("H68").Value = 15997.5
$ws.Range("I68").Value = 10000
$ws.Range("J68").Value = 17996.666
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 17996.666
$ws.Range("M68").Value = -9251
$ws.Range("N68").Value = -19494.666
$ws.Range("H71").Value = 15997.5
$ws.Range("I71").Value = 10000
$ws.Range("J71").Value = 17996.666
$ws.Range("K71").Value = 30000
$ws.Range("L71").Value = 53989.99800000001
$ws.Range("M71").Value = -26256
$ws.Range("N71").Value = -61477.99800000001
$ws.Range("H74").Value = 16838.334
$ws.Range("J74").Value = 18606
$ws.Range("L74").Value = 18606
$ws.Range("N74").Value = -20354
$ws.Range("H77").Value = 16838.334
$ws.Range("J77").Value = 18606
$ws.Range("L77").Value = 55818
$ws.Range("N77").Value = -64554
$ws.Range("H98").Value = 45065.715
$ws.Range("J98").Value = 45065.715
$ws.Range("L98").Value = 45065.715
$ws.Range("N98").Value = -49557.715
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -38180
$ws.Range("H134").Value = 921.45
$ws.Range("I134").Value = 893.79486
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2681.38458
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -146.3845799999999
$ws.Range("N134").Value = -11070
$ws.Range("H136").Value = 3300.1904
$ws.Range("I136").Value = 771.7143
$ws.Range("J136").Value = 8357.143
$ws.Range("K136").Value = 2315.1429
$ws.Range("L136").Value = 25071.429
$ws.Range("M136").Value = 234.8571000000002
$ws.Range("N136").Value = -30171.429
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 89286540
$ws.Range("J131").Value = 125000890
$ws.Range("L131").Value = 375002670
$ws.Range("N131").Value = -375012750
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 950
$ws.Range("J13").Value = 950
$ws.Range("L13").Value = 950
$ws.Range("N13").Value = -1228
$ws.Range("H126").Value = 1295.3914
$ws.Range("J126").Value = 1511.4286
$ws.Range("L126").Value = 4534.2858
$ws.Range("N126").Value = -9474.2858
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4421
$ws.Range("I7").Value = 4421
$ws.Range("K7").Value = 4421
$ws.Range("M7").Value = -4309
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H122").Value = 2659.9375
$ws.Range("I122").Value = 2109.9
$ws.Range("J122").Value = 3576.6667
$ws.Range("K122").Value = 6329.700000000001
$ws.Range("L122").Value = 10730.0001
$ws.Range("M122").Value = -3879.700000000001
$ws.Range("N122").Value = -15630.0001
$ws.Range("H126").Value = 4421
$ws.Range("I126").Value = 4421
$ws.Range("K126").Value = 13263
$ws.Range("M126").Value = -10793
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7459.154
$ws.Range("I122").Value = 5760.0625
$ws.Range("J122").Value = 10177.7
$ws.Range("K122").Value = 17280.1875
$ws.Range("L122").Value = 30533.1
$ws.Range("M122").Value = -14830.1875
$ws.Range("N122").Value = -35433.10000000001
$ws.Range("H123").Value = 51795
$ws.Range("J123").Value = 51795
$ws.Range("L123").Value = 51795
$ws.Range("N123").Value = -61595
$ws.Range("H140").Value = 49300
$ws.Range("J140").Value = 49300
$ws.Range("L140").Value = 49300
$ws.Range("N140").Value = -59660
